$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45228, "10:00", 2090, "amazon", "preto"),
    @(45229, "20:58", 2090, "amazon", "preto"),
    @(45230, "12:41", 2090, "amazon", "preto")
)

$row = 32
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row++
}
